$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ativação date: 01/01/2020 -> 01/01/2023
# Leading apostrophe forces this to be stored as literal text instead of
# being auto-converted to a date serial value.
$ws.Range("B8").Value = "'01/01/2023"
$ws.Range("C8").Value = "'01/01/2023"

# New professor (5840712 - Ângelo Capri Neto) added to the course.
# He shows up both on the "Objetivos:" row (B10/C10) and the
# "Programa resumido:" row (B13/C13).
$ws.Range("B10").Value = "5840712 - Ângelo Capri Neto"
$ws.Range("C10").Value = "5840712 - Ângelo Capri Neto"
$ws.Range("B13").Value = "5840712 - Ângelo Capri Neto"
$ws.Range("C13").Value = "5840712 - Ângelo Capri Neto"

# Critério text (row 19)
$ws.Range("B19").Value = "A avaliação tem como requisito quantificar as competências adquiridas conforme objetivadas.Duas provas escritas (P1 e P2) e listas de exercícios de acompanhamento continuado. A partir das notas das listas de exercício será calculada a média, LE."
$ws.Range("C19").Value = "A avaliação tem como requisito quantificar as competências adquiridas conforme objetivadas.Duas provas escritas (P1 e P2) e listas de exercícios de acompanhamento continuado. A partir das notas das listas de exercício será calculada a média, LE."

# Norma de recuperação value (row 20)
$ws.Range("B20").Value = "NF = (P1 + P2 + LE) /3"
$ws.Range("C20").Value = "NF = (P1 + P2 + LE) /3"

# Bibliografia row value (row 21) - recovery exam description
$ws.Range("B21").Value = "Será realizada uma prova escrita valendo de zero a dez (NR) e a média final calculada pela equação: NF + NR"
$ws.Range("C21").Value = "Será realizada uma prova escrita valendo de zero a dez (NR) e a média final calculada pela equação: NF + NR"
